$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-01 Friday", "2024-11-02 Saturday"),
    @("32×64=", "73×16="),
    @("58×91=", "17×61="),
    @("27×57=", "90×29="),
    @("47×48=", "24×83="),
    @("40×17=", "29×83="),
    @("37×83=", "37×41="),
    @("92×80=", "57×52="),
    @("50×24=", "76×68="),
    @("79×66=", "46×79="),
    @("51×76=", "52×54="),
    @("13×30=", "87×63="),
    @("30×72=", "83×60="),
    @("31×73=", "34×87="),
    @("96×77=", "93×41="),
    @("66×55=", "24×56="),
    @("75×35=", "29×83="),
    @("77×20=", "38×96="),
    @("62×84=", "89×27="),
    @("93×66=", "63×59="),
    @("14×47=", "55×19="),
    @("66×51=", "81×98="),
    @("74×92=", "87×77="),
    @("48×78=", "81×19="),
    @("91×21=", "32×93="),
    @("78×97=", "40×26=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
